$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: becomes the "Spillkråka" (Dryocopus martius) record ---
$ws.Range("A10").Value = 111964494
$ws.Range("B10").Value = 56446
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 100049
$ws.Range("F10").Value = "Spillkråka"
$ws.Range("G10").Value = "Dryocopus martius"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = "adult"
$ws.Range("M10").Value = "förbiflygande"
$ws.Range("P10").Value = "Stenstorp SSO 1470 m, Ög"
$ws.Range("Q10").Value = 575346
$ws.Range("R10").Value = 6509958
$ws.Range("Z10").Value = "10:30"
$ws.Range("AB10").Value = "10:30"
$ws.Range("AF10").Value = ""
$ws.Range("AI10").Value = "Äldre barrskog"

# --- Row 11: becomes the "Blåmossa" (Leucobryum glaucum) record ---
$ws.Range("A11").Value = 111964621
$ws.Range("B11").Value = 93553
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 2180
$ws.Range("F11").Value = "Blåmossa"
$ws.Range("G11").Value = "Leucobryum glaucum"
$ws.Range("H11").Value = "(Hedw.) Ångstr."
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("M11").Value = ""
$ws.Range("P11").Value = "Stenstorp SSO 1660 m, Ög"
$ws.Range("Q11").Value = 575609
$ws.Range("R11").Value = 6509825
$ws.Range("Z11").Value = ""
$ws.Range("AB11").Value = ""
$ws.Range("AF11").Value = ""
$ws.Range("AI11").Value = "Barrskog"

# --- Row 12: small taxon-order id tweak ---
$ws.Range("B12").Value = 103781
